$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster changes from "Resolving-Mac" to "Inflammatory-Mac"
$ws.Range("A2").Value = "Inflammatory-Mac"

# Row 3: Target cluster changes from "MuSCs" to "ECs"
$ws.Range("D3").Value = "ECs"

# Row 2 numeric updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.119110666666667
$ws.Range("H2").Value = 3.357332
$ws.Range("I2").Value = 0.5726510027906514
$ws.Range("J2").Value = 0.5726510027906513
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0987018304311111
$ws.Range("R2").Value = 0.8883164738799999
$ws.Range("S2").Value = 0.5726510027906514
$ws.Range("T2").Value = 0.5726510027906513

# Row 3 numeric updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8351523333333333
$ws.Range("H3").Value = 2.505457
$ws.Range("I3").Value = 0.4273489972093487
$ws.Range("J3").Value = 0.4273489972093487
$ws.Range("M3").Value = 0.08819666666666666
$ws.Range("N3").Value = 0.26459
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.07365765195888888
$ws.Range("R3").Value = 0.66291886763
$ws.Range("S3").Value = 0.4273489972093487
$ws.Range("T3").Value = 0.4273489972093487
